# Add a new "legume fodder" land-use block.
# The existing data is laid out in repeating 3-column groups (one group per
# land use), for rows 1-10. A new group is inserted right after the group
# currently occupying columns S:U, by inserting 3 blank columns at V:X and
# filling them with a copy of the S:U group (columns before it stay the
# same, since no new rotation data exists yet - see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new (blank) columns before the current column V. This shifts the
# old V:BB data right to Y:BE, growing the used range to A1:BE10.
$ws.Columns("V:X").Insert()

# Populate the newly inserted columns with a copy of the S:U block so the
# new land-use group mirrors the one it was duplicated from.
$ws.Range("S1:U10").Copy($ws.Range("V1"))
